{"js": "// Find the existing bullet after which the three new bullets must be inserted.\nconst body = context.document.body;\nconst anchorText =\n  \"\u2022 Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\";\n\nconst results = body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const newBullets = [\n    \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n    \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n    \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n  ];\n\n  // Insert each bullet as its own new paragraph, right after the anchor\n  // paragraph (and after each previously-inserted one, preserving order).\n  let insertionPoint = results.items[0];\n  for (const bulletText of newBullets) {\n    insertionPoint = insertionPoint.insertParagraph(bulletText, Word.InsertLocation.after);\n  }\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the existing bullet that the new bullets must be inserted after.\n$anchorText = \"Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\"\n$rng = $d.Content\n$found = $rng.Find.Execute($anchorText)\n\nif ($found) {\n    $bullet1 = \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\"\n    $bullet2 = \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\"\n    $bullet3 = \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n\n    # $rng now spans the found text; inserting after it (with a leading\n    # paragraph mark) adds three new paragraphs right after the anchor\n    # paragraph, before the next existing bullet.\n    $rng.InsertAfter(\"`r\" + $bullet1 + \"`r\" + $bullet2 + \"`r\" + $bullet3)\n}\n"}
